$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numbers as plain text (e.g. "41.853.11", "0.0658")
# rather than as real numbers. Writing a numeric-looking string via .Value would
# normally let Excel auto-convert it to a Number; force Text first, then restore the
# original (default/"Normal") cell style so no extra formatting is introduced.
$dCells = @("D2","D3","D5","D7","D9","D10","D13","D14","D15","D18","D19","D21","D22","D23","D24","D27","D28","D32","D34","D36","D41","D42","D45","D46","D47","D48","D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.853.11"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "2.231.89"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "232.72"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("D7").Value = "60.77"
$ws.Range("E7").Value = "  -5.52%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.406"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "58.07"
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "2.562.01"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "15.69"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "22.71"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "2.246.92"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "41.737.57"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "72.54"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "247.55"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "9.63"
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("D28").Value = "169.17"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").Value = "2.64"
$ws.Range("E32").Value = "  -6.52%  "
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  +5.85%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "0.0658"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("E37").Value = "  -8.15%  "
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  -5.36%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").Value = "0.000236"
$ws.Range("E41").Value = "  +13.05%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0241"
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "4.48"
$ws.Range("E45").Value = "  -9.10%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "98.69"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0959"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").Value = "1.470.84"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").Value = "16.61"
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("E50").Value = "  +7.93%  "
$ws.Range("E51").Value = "  -2.73%  "

# Restore default styling on the Price cells (remove the temporary Text format)
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
